{"js": "// Fix two typos in the \"Domn\u00edv\u00e1m se, \u017ee, aby se mohl akcion\u00e1\u0159...\" paragraph:\n//   1) \"z\u00e1chov\u00e1no\" -> \"zachov\u00e1no\"   (wrong vowel length on the first syllable)\n//   2) \"akcion\u00e1\u0159u\"  -> \"akcion\u00e1\u0159\u016f\"  (wrong case ending: genitive plural)\n//\n// Both strings are unique in the document, so a direct search/replace is safe\n// and unambiguous.\n\nconst body = context.document.body;\n\n// 1) z\u00e1chov\u00e1no -> zachov\u00e1no\nconst hits1 = body.search(\"z\u00e1chov\u00e1no\", { matchCase: true, matchWholeWord: false });\nhits1.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < hits1.items.length; i++) {\n  hits1.items[i].insertText(\"zachov\u00e1no\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) akcion\u00e1\u0159u -> akcion\u00e1\u0159\u016f\nconst hits2 = body.search(\"akcion\u00e1\u0159u\", { matchCase: true, matchWholeWord: false });\nhits2.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < hits2.items.length; i++) {\n  hits2.items[i].insertText(\"akcion\u00e1\u0159\u016f\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix two typos in the \"Domn\u00edv\u00e1m se, \u017ee, aby se mohl akcion\u00e1\u0159...\" paragraph\n# (BodyText style):\n#   1) \"z\u00e1chov\u00e1no\" -> \"zachov\u00e1no\"   (wrong vowel length on the first syllable)\n#   2) \"akcion\u00e1\u0159u\"  -> \"akcion\u00e1\u0159\u016f\"  (wrong case ending: genitive plural)\n#\n# Both strings are unique in the document, so Find/Replace against the whole\n# document content is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n# 1) z\u00e1chov\u00e1no -> zachov\u00e1no\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"z\u00e1chov\u00e1no\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"zachov\u00e1no\"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# 2) akcion\u00e1\u0159u -> akcion\u00e1\u0159\u016f\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"akcion\u00e1\u0159u\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"akcion\u00e1\u0159\u016f\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
